# "changed news limits (summary limit is now 1000 chars) and checks and added new tests and questions"
#
# Adds 7 new "News" related test cases (rows 25-31) into the Tabelle1 test
# sheet, widens the last two table columns a bit, grows the trailing filler
# rows by two (to 402) and resizes the "Tabelle1" Excel table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ------------------------------------------------------
# Columns E:F (Eingaben / Erwartetes Ergebnis) get a little wider.
# Columns C:D (Kurzbeschreibung / Vorbedingungen) stay as they are.
$ws.Columns.Item(5).ColumnWidth = 34.9
$ws.Columns.Item(6).ColumnWidth = 34.9

# --- Row 25 --------------------------------------------------------------
$ws.Range("A25").Value = 16
$ws.Range("B25").Value = "NF"
$ws.Range("C25").Value = "alle News anzeigen"
$ws.Range("D25").Value = "es gibt bereits mehr als 20 News Einträge im System"
$ws.Range("E25").Value = 'einen News Eintrag öffnen - danach wieder schließen, danach unten auf "Gelesene Anzeigen?" klicken'
$ws.Range("F25").Value = "der gerade gelesene Eintrag sollte auch zu sehen sein und beim nach unten scrollen sollten gepaged alle News geladen werden"
$ws.Rows.Item(25).RowHeight = 60

# --- Row 26 --------------------------------------------------------------
$ws.Range("A26").Value = 17
$ws.Range("B26").Value = "NF"
$ws.Range("C26").Value = "alle ungelesenen News anzeigen"
$ws.Range("D26").Value = "es gibt bereits mehr als 20 News Einträge im System"
$ws.Range("E26").Value = "einen News Eintrag öffnen und danach wieder schließen"
$ws.Range("F26").Value = "der gerade gelesene Eintrag ist jetzt nicht mehr in der Liste, beim nach unten scrollen werden die nächsten ungelesenen News gepaged geladen"
$ws.Rows.Item(26).RowHeight = 60

# --- Row 27 --------------------------------------------------------------
$ws.Range("A27").Value = 18
$ws.Range("B27").Value = "NF"
$ws.Range("C27").Value = "neuen News Beitrag ohne Bild erstellen"
$ws.Range("D27").Value = "/"
$ws.Range("E27").Value = 'Klick auf "News hinzufügen" und alle Felder bis auf das Bild ausfüllen und mit News veröffentlichen bestätigen'
$ws.Range("F27").Value = "das Fenster schließt sich und der Beitrag wird ganz oben angezeigt (als ungelesener Beitrag)"
$ws.Rows.Item(27).RowHeight = 60

# --- Row 28 --------------------------------------------------------------
$ws.Range("A28").Value = 19
$ws.Range("B28").Value = "FF"
$ws.Range("C28").Value = "neuen unvollständigen News Beitrag erstellen"
$ws.Range("D28").Value = "/"
$ws.Range("E28").Value = 'Klick auf "News hinzufügen" und entweder Titel, Zusammenfassung oder Text nicht ausfüllen und mit News veröffentlichen bestätigen'
$ws.Range("F28").Value = "es kommt eine Fehlermeldung die (mehrsprachig) anzeigt, welches Feld nicht ausgefüllt wurde"
$ws.Rows.Item(28).RowHeight = 60

# --- Row 29 --------------------------------------------------------------
$ws.Range("A29").Value = 20
$ws.Range("B29").Value = "FF"
$ws.Range("C29").Value = "neuen News Beitrag erstellen mit zu langen Texten"
$ws.Range("D29").Value = "/"
$ws.Range("E29").Value = 'Klick auf "News hinzufügen" und entweder beim Titel mehr als 100 Zeichen, bei der Zusammenfassung mehr als 1.000 oder beim Text mehr als 10.000 Zeichen eintragen und mit News veröffentlichen bestätigen'
$ws.Range("F29").Value = "es kommt eine Fehlermeldung die (mehrsprachig) anzeigt, welcher Text zu lang ist"
$ws.Rows.Item(29).RowHeight = 90

# --- Row 30 --------------------------------------------------------------
$ws.Range("A30").Value = 21
$ws.Range("B30").Value = "NF"
$ws.Range("C30").Value = "neuen News Beitrag mit Bild erstellen"
$ws.Range("D30").Value = "/"
$ws.Range("E30").Value = 'Klick auf "News hinzufügen" und alle Felder inklusive einem Bild (kleiner als 5 MB) ausfüllen und mit News veröffentlichen bestätigen'
$ws.Range("F30").Value = "das Fenster schließt sich und der Beitrag wird ganz oben angezeigt (als ungelesener Beitrag) - wenn man diesen nun öffnet sieht man auch das Bild"
$ws.Rows.Item(30).RowHeight = 75

# --- Row 31 --------------------------------------------------------------
$ws.Range("A31").Value = 22
$ws.Range("B31").Value = "FF"
$ws.Range("C31").Value = "neuen News Beitrag mit einem Bild > 5 MB erstellen"
$ws.Range("D31").Value = "/"
$ws.Range("E31").Value = 'Klick auf "News hinzufügen" und alle Felder inklusive einem Bild größer als 5 MB) ausfüllen und mit News veröffentlichen bestätigen'
$ws.Range("F31").Value = "es kommt eine Fehlermeldung die (mehrsprachig) anzeigt, dass das Bild zu groß ist"
$ws.Rows.Item(31).RowHeight = 60

# --- Extend the trailing filler rows (table grew by two data rows, so the
#     blank spacer rows after it now end at 402 instead of 400) -----------
$ws.Range("A401").HorizontalAlignment = -4108
$ws.Range("A401").VerticalAlignment = -4108
$ws.Rows.Item(401).RowHeight = 20.1

$ws.Range("A402").HorizontalAlignment = -4108
$ws.Range("A402").VerticalAlignment = -4108
$ws.Rows.Item(402).RowHeight = 20.1

# --- Grow the "Tabelle1" table / autofilter range to match -----------------
$ws.ListObjects.Item(1).Resize($ws.Range("A9:F402"))

# --- Selection / view state (best-effort; matches the saved selection) ----
$ws.Range("H27").Select() | Out-Null

Write-Output "applied edits"
